# Added time for group work on 22-02-2016
# Record 4 hours (4/24 of a day) of "gemeenschappelijk" (group) time in week 2 (column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 4 / 24

# Move the active selection to the cell below the edited one, matching
# where the cursor ends up after typing a value into E10 and pressing Enter.
[void]$ws.Range("E11").Select()
